# Update the "Förändrad" (Changed) date column (C) for all data rows
# from 45439 (2024-05-27) to 45440 (2024-05-28).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 28; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45439) {
        $cell.Value = 45440
    }
}
